$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 9197.333000000001
$ws.Range("I64").Value = 3246
$ws.Range("K64").Value = 3246
$ws.Range("M64").Value = -2998
$ws.Range("H67").Value = 9197.333000000001
$ws.Range("I67").Value = 3246
$ws.Range("K67").Value = 3246
$ws.Range("M67").Value = -2388
$ws.Range("H87").Value = 67666
$ws.Range("J87").Value = 67666
$ws.Range("L87").Value = 67666
$ws.Range("N87").Value = -70162
$ws.Range("H90").Value = 67666
$ws.Range("J90").Value = 67666
$ws.Range("L90").Value = 202998
$ws.Range("N90").Value = -215478
$ws.Range("H92").Value = 26968.254
$ws.Range("I92").Value = 13226.904
$ws.Range("J92").Value = 71363.38
$ws.Range("K92").Value = 13226.904
$ws.Range("L92").Value = 71363.38
$ws.Range("M92").Value = -11978.904
$ws.Range("N92").Value = -73859.38
$ws.Range("H106").Value = 3219.3572
$ws.Range("I106").Value = 3357.1
$ws.Range("J106").Value = 2875
$ws.Range("K106").Value = 3357.1
$ws.Range("L106").Value = 2875
$ws.Range("M106").Value = -2726.1
$ws.Range("N106").Value = -4137
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 754
$ws.Range("N113").ClearContents()
$ws.Range("H127").Value = 1408.8572
$ws.Range("I127").Value = 1408.8572
$ws.Range("K127").Value = 4226.571599999999
$ws.Range("M127").Value = 733.4284000000007
$ws.Range("H128").Value = 66890
$ws.Range("J128").Value = 66890
$ws.Range("L128").Value = 66890
$ws.Range("N128").Value = -76850

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 48666
$ws.Range("J44").Value = 38399.2
$ws.Range("L44").Value = 38399.2
$ws.Range("N44").Value = -39375.2
$ws.Range("H63").Value = 2160.5715
$ws.Range("I63").Value = 2160.5715
$ws.Range("K63").Value = 2160.5715
$ws.Range("M63").Value = -1474.5715
$ws.Range("H66").Value = 2160.5715
$ws.Range("I66").Value = 2160.5715
$ws.Range("K66").Value = 10802.8575
$ws.Range("M66").Value = -7370.8575
$ws.Range("H122").Value = 100997.22
$ws.Range("I122").Value = 999.6667
$ws.Range("J122").Value = 150996
$ws.Range("K122").Value = 2999.0001
$ws.Range("L122").Value = 452988
$ws.Range("M122").Value = -549.0001000000002
$ws.Range("N122").Value = -457888

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17945.824
$ws.Range("I82").Value = 3426.4
$ws.Range("J82").Value = 38687.855
$ws.Range("K82").Value = 3426.4
$ws.Range("L82").Value = 38687.855
$ws.Range("M82").Value = -3043.4
$ws.Range("N82").Value = -39453.855
$ws.Range("H85").Value = 17945.824
$ws.Range("I85").Value = 3426.4
$ws.Range("J85").Value = 38687.855
$ws.Range("K85").Value = 3426.4
$ws.Range("L85").Value = 38687.855
$ws.Range("M85").Value = -2100.4
$ws.Range("N85").Value = -41339.855
$ws.Range("H86").Value = 2086.647
$ws.Range("I86").Value = 1323.375
$ws.Range("J86").Value = 2765.111
$ws.Range("K86").Value = 1323.375
$ws.Range("L86").Value = 2765.111
$ws.Range("M86").Value = -200.375
$ws.Range("N86").Value = -5011.111
$ws.Range("H89").Value = 2086.647
$ws.Range("I89").Value = 1323.375
$ws.Range("J89").Value = 2765.111
$ws.Range("K89").Value = 6616.875
$ws.Range("L89").Value = 13825.555
$ws.Range("M89").Value = -1000.875
$ws.Range("N89").Value = -25057.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 11832.833
$ws.Range("I36").Value = 19499.5
$ws.Range("J36").Value = 7999.5
$ws.Range("K36").Value = 19499.5
$ws.Range("L36").Value = 7999.5
$ws.Range("M36").Value = -19111.5
$ws.Range("N36").Value = -8775.5
$ws.Range("H40").Value = 11832.833
$ws.Range("I40").Value = 19499.5
$ws.Range("J40").Value = 7999.5
$ws.Range("K40").Value = 19499.5
$ws.Range("L40").Value = 7999.5
$ws.Range("M40").Value = -19339.5
$ws.Range("N40").Value = -8319.5
$ws.Range("H41").Value = 38087.332
$ws.Range("J41").Value = 39098.25
$ws.Range("L41").Value = 39098.25
$ws.Range("N41").Value = -39954.25
$ws.Range("H62").Value = 45185.81
$ws.Range("I62").Value = 84658.08
$ws.Range("J62").Value = 5713.5386
$ws.Range("K62").Value = 84658.08
$ws.Range("L62").Value = 5713.5386
$ws.Range("M62").Value = -84034.08
$ws.Range("N62").Value = -6961.5386
$ws.Range("H65").Value = 45185.81
$ws.Range("I65").Value = 84658.08
$ws.Range("J65").Value = 5713.5386
$ws.Range("K65").Value = 423290.4
$ws.Range("L65").Value = 28567.693
$ws.Range("M65").Value = -420170.4
$ws.Range("N65").Value = -34807.693
$ws.Range("H70").Value = 59999
$ws.Range("J70").Value = 59999
$ws.Range("L70").Value = 59999
$ws.Range("N70").Value = -60629
$ws.Range("H73").Value = 59999
$ws.Range("J73").Value = 59999
$ws.Range("L73").Value = 59999
$ws.Range("N73").Value = -62183

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 123830.5
$ws.Range("J37").Value = 123830.5
$ws.Range("L37").Value = 371491.5
$ws.Range("N37").Value = -371715.5
$ws.Range("H70").Value = 3752.75
$ws.Range("J70").Value = 4000
$ws.Range("L70").Value = 12000
$ws.Range("N70").Value = -12630
$ws.Range("H73").Value = 3752.75
$ws.Range("J73").Value = 4000
$ws.Range("L73").Value = 12000
$ws.Range("N73").Value = -14184
$ws.Range("H132").Value = 2306.6956
$ws.Range("I132").Value = 1384.2142
$ws.Range("J132").Value = 3741.6667
$ws.Range("K132").Value = 12457.9278
$ws.Range("L132").Value = 33675.0003
$ws.Range("M132").Value = -9927.927799999999
$ws.Range("N132").Value = -38735.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 56672.15
$ws.Range("I102").Value = 112807
$ws.Range("J102").Value = 10743.637
$ws.Range("K102").Value = 112807
$ws.Range("L102").Value = 10743.637
$ws.Range("M102").Value = -111185
$ws.Range("N102").Value = -13987.637

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5085.8
$ws.Range("I46").Value = 868.2
$ws.Range("K46").Value = 868.2
$ws.Range("M46").Value = -680.2
